# This script applies the data refresh for the cryptos worksheet.
# It updates Price (column D) and Volume(1h) (column E) values for
# the listed coins, and reflects the reordering of the ImmutableX /
# EthereumClassic rows (34 and 35) along with their refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.891.12"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "'2.737.17"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'565.99"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "'161.31"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "'5.63"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "'0.377"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "'3.219.88"
$ws.Range("D14").Value = "'26.99"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "'63.667.27"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "'2.739.62"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'12.35"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'4.75"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").Value = "'356.11"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D23").Value = "'0.521"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").Value = "'64.21"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'8.38"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'0.0₃0913"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'2.00"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").Value = "'1.36"
$ws.Range("E30").Value = "  +9.16%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("D32").Value = "'167.12"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "'4.94"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'20.08"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").Value = "'0.979"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "'349.22"
$ws.Range("E39").Value = "  +5.33%  "
$ws.Range("D40").Value = "'6.30"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").Value = "'4.08"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "'38.63"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "'21.86"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("D44").Value = "'21.00"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "'0.0584"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").Value = "'0.632"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'0.0251"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").Value = "'0.0996"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "'132.45"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "'11.07"
$ws.Range("E51").Value = "  +0.40%  "

Write-Host "Updated cryptos list"
